$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 53363
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 6846
